$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OnboardingRegisteredInstitute")

# ---------------------------------------------------------------------------
# 1. Drop the old "TC Name" column (B). Role/username/password shift left
#    into B/C/D.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).Delete()

# ---------------------------------------------------------------------------
# 2. Insert 9 new columns (E..M) to hold the faculty-creation test data.
# ---------------------------------------------------------------------------
$ws.Range("E1:M1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 3. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.3
$ws.Columns.Item(2).ColumnWidth = 18.49
$ws.Columns.Item(3).ColumnWidth = 31.57
$ws.Columns.Item(4).ColumnWidth = 20.72
$ws.Columns.Item(5).ColumnWidth = 14.46
$ws.Columns.Item(6).ColumnWidth = 13.9
$ws.Columns.Item(7).ColumnWidth = 16.69
$ws.Columns.Item(8).ColumnWidth = 16.83
$ws.Columns.Item(9).ColumnWidth = 11.96
$ws.Columns.Item(10).ColumnWidth = 14.18
$ws.Columns.Item(11).ColumnWidth = 18.49
$ws.Columns.Item(12).ColumnWidth = 11.12
$ws.Columns.Item(13).ColumnWidth = 26.28

# ---------------------------------------------------------------------------
# 4. Header row (row 1) text for the new columns.
# ---------------------------------------------------------------------------
$ws.Range("E1").Value2 = "Username"
$ws.Range("F1").Value2 = "ChooseSubject"
$ws.Range("G1").Value2 = "Email"
$ws.Range("H1").Value2 = "PhoneNumber"
$ws.Range("I1").Value2 = "Address"
$ws.Range("J1").Value2 = "Biography"
$ws.Range("K1").Value2 = "Occupation"
$ws.Range("L1").Value2 = "Password"
$ws.Range("M1").Value2 = "ConfirmPassword"

# Header formatting: yellow fill + thin border all round, matching the
# existing TC No/Role/username/password header cells. The new header cells
# are not bold (unlike the original ones).
$hdrNew = $ws.Range("E1:M1")
$hdrNew.Font.Name = "Calibri"
$hdrNew.Font.Size = 11
$hdrNew.Font.Bold = $false
$hdrNew.Font.ColorIndex = 1
$hdrNew.Interior.Color = 65535
$hdrNew.Borders.LineStyle = 1

# Re-affirm the original header cells keep their bold/yellow/border look.
$hdrOld = $ws.Range("A1:D1")
$hdrOld.Font.Name = "Calibri"
$hdrOld.Font.Size = 11
$hdrOld.Font.Bold = $true
$hdrOld.Font.ColorIndex = 1
$hdrOld.Interior.Color = 65535
$hdrOld.Borders.LineStyle = 1

$ws.Rows.Item(1).RowHeight = 13.8

# ---------------------------------------------------------------------------
# 5. Data row (row 2) values.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value2 = "faculty123"
$ws.Range("F2").Value2 = "Biology"
$ws.Range("G2").Value2 = "testFaculty@gmail.com"
$ws.Range("H2").Value2 = 8956895689
$ws.Range("I2").Value2 = "bangalore"
$ws.Range("J2").Value2 = "biography"
$ws.Range("K2").Value2 = "occupation"
$ws.Range("L2").Value2 = "password"
$ws.Range("M2").Value2 = "password"

# Data formatting: plain (no fill) + thin border all round for the whole row.
$body = $ws.Range("A2:M2")
$body.Font.Name = "Calibri"
$body.Font.Size = 11
$body.Font.Bold = $false
$body.Font.ColorIndex = 1
$body.Interior.ColorIndex = 0
$body.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 6. Hyperlinks. Delete the stale hyperlink (it stayed anchored on the old
#    D2 cell) and recreate the institute-email link on its new home (C2),
#    plus add the new faculty-email hyperlink on G2.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:sampleTestInstitute@gmail.com", "", "", "sampleTestInstitute@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:testFaculty@gmail.com", "", "", "testFaculty@gmail.com")

# Undo the automatic "Hyperlink" style (blue/underlined) the Add() calls
# above apply, restoring the plain bordered look used by the rest of row 2.
# (Done per-cell rather than via a union range - multi-area ranges don't
# reliably broadcast writes in this host.)
foreach ($addr in "C2", "G2") {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Underline = $false
    $cell.Font.ColorIndex = 1
    $cell.Interior.ColorIndex = 0
    $cell.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# 7. Selection matches the committed sheet (active cell D2).
# ---------------------------------------------------------------------------
$ws.Range("D2").Select()
